# Add task "UserShippingAddress" (task #11) and finish task #10 ("BookDetail")
# on both the Back-end and Front-end sheets.

$wb = $excel.ActiveWorkbook

$backend = $wb.Worksheets.Item("Back-end")
$frontend = $wb.Worksheets.Item("Front-end")

# --- Back-end sheet -------------------------------------------------
# Order of writes matters so that new shared strings are created in the
# same order they appear in the target workbook: 4/5, 3h, 7/5,
# UserShippingAddress, 8/5, list...

# Row 16 (task #10 "BookDetail") - fill in remaining Actual/Status info
$backend.Range("I16").Value = "4/5"
$backend.Range("D16").Value = "3h"
$backend.Range("G16").Value = "3h"
$backend.Range("H16").Value = "3/5"
$backend.Range("J16").Value = "Done"

# Row 17 (task #11 "UserShippingAddress") - brand new task row
$backend.Range("E17").Value = "7/5"
$backend.Range("C17").Value = "UserShippingAddress"
$backend.Range("F17").Value = "8/5"
$backend.Range("B17").Value = "bookstore"
$backend.Range("D17").Value = "8h"
$backend.Range("G17").Value = "8h"
$backend.Range("H17").Value = "7/5"
$backend.Range("I17").Value = "8/5"
$backend.Range("J17").Value = "Done"
$backend.Range("K17").Value = "list, add, update, delete, setDefault"
$backend.Range("K17").WrapText = $true

# --- Front-end sheet -------------------------------------------------
# Row 16 (task #10 "BookDetail") - fill in remaining Actual/Status info
$frontend.Range("G16").Value = "1h"
$frontend.Range("H16").Value = "3/5"
$frontend.Range("I16").Value = "4/5"
$frontend.Range("J16").Value = "Done"

# Row 17 (task #11 "UserShippingAddress") - brand new task row
$frontend.Range("B17").Value = "bookstore"
$frontend.Range("C17").Value = "UserShippingAddress"
$frontend.Range("D17").Value = "8h"
$frontend.Range("E17").Value = "7/5"
$frontend.Range("F17").Value = "8/5"
$frontend.Range("G17").Value = "6h"
$frontend.Range("H17").Value = "7/5"
$frontend.Range("I17").Value = "8/5"
$frontend.Range("J17").Value = "Done"

# --- Selection state (cosmetic, matches authored workbook) -----------
$backend.Range("K26").Select() | Out-Null
$frontend.Range("E17").Select() | Out-Null
$backend.Activate() | Out-Null
